$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 7 (row 19) "RR vs DC" - enter points scored by each of the 9 players
$ws.Range("E19").Value = 30
$ws.Range("H19").Value = 20
$ws.Range("K19").Value = 50
$ws.Range("N19").Value = 80
$ws.Range("Q19").Value = 100
$ws.Range("T19").Value = 0
$ws.Range("W19").Value = 70
$ws.Range("Z19").Value = 40
$ws.Range("AC19").Value = 60
